# Refactor synthetic array /3
# Rename the "noir" color entry (and its emoji swatch) to "bleu", reusing
# a new set of book emoji for the four status swatches while keeping the
# existing rouge/orange/vert entries untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange

# Map of old shared-string text -> new shared-string text.
$replacements = @{
    "⬛"   = "📘"
    "🟥"   = "📕"
    "🟧"   = "📙"
    "🟩"   = "📗"
    "noir" = "bleu"
}

foreach ($old in $replacements.Keys) {
    $new = $replacements[$old]
    $used.Replace($old, $new, 1, 1, $false, $false) | Out-Null
}
